$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four trailing rows (ma.L2, ma.S.L12, ma.S.L24, sigma2 old rows)
# that are no longer present after the outlier-detection re-run.
$ws.Rows("17:20").Delete()

# Update rows 2-16 with the new coefficients from the AO-only outlier run.
$ws.Range("A2").Value = "AO2008Sep"
$ws.Range("B2").Value = -39880
$ws.Range("C2").Value = 95000
$ws.Range("D2").Value = -0.42
$ws.Range("E2").Value = 0.675
$ws.Range("F2").Value = -226000
$ws.Range("G2").Value = 146000
$ws.Range("A3").Value = "AO2009May"
$ws.Range("B3").Value = -64500
$ws.Range("C3").Value = 96900
$ws.Range("D3").Value = -0.666
$ws.Range("E3").Value = 0.506
$ws.Range("F3").Value = -254000
$ws.Range("G3").Value = 125000
$ws.Range("A4").Value = "AO2017Aug"
$ws.Range("B4").Value = -72080
$ws.Range("C4").Value = 91000
$ws.Range("D4").Value = -0.792
$ws.Range("E4").Value = 0.429
$ws.Range("F4").Value = -251000
$ws.Range("G4").Value = 106000
$ws.Range("A5").Value = "AO2020Apr"
$ws.Range("B5").Value = -318200
$ws.Range("C5").Value = 29600
$ws.Range("D5").Value = -10.76
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = -376000
$ws.Range("G5").Value = -260000
$ws.Range("A6").Value = "AO2020May"
$ws.Range("B6").Value = -377700
$ws.Range("C6").Value = 24200
$ws.Range("D6").Value = -15.639
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = -425000
$ws.Range("G6").Value = -330000
$ws.Range("A7").Value = "AO2020Aug"
$ws.Range("B7").Value = -324900
$ws.Range("C7").Value = 29400
$ws.Range("D7").Value = -11.039
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = -383000
$ws.Range("G7").Value = -267000
$ws.Range("A8").Value = "ar.L1"
$ws.Range("B8").Value = -0.5704
$ws.Range("C8").Value = 0.134
$ws.Range("D8").Value = -4.269
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = -0.832
$ws.Range("G8").Value = -0.308
$ws.Range("A9").Value = "ar.L2"
$ws.Range("B9").Value = 0.2459
$ws.Range("C9").Value = 0.076
$ws.Range("D9").Value = 3.224
$ws.Range("E9").Value = 0.001
$ws.Range("F9").Value = 0.096
$ws.Range("G9").Value = 0.395
$ws.Range("A10").Value = "ar.L3"
$ws.Range("B10").Value = 0.5436
$ws.Range("C10").Value = 0.103
$ws.Range("D10").Value = 5.26
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0.341
$ws.Range("G10").Value = 0.746
$ws.Range("A11").Value = "ar.L4"
$ws.Range("B11").Value = 0.5141
$ws.Range("C11").Value = 0.082
$ws.Range("D11").Value = 6.267
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0.353
$ws.Range("G11").Value = 0.675
$ws.Range("A12").Value = "ma.L1"
$ws.Range("B12").Value = 1.2333
$ws.Range("C12").Value = 0.151
$ws.Range("D12").Value = 8.15
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0.937
$ws.Range("G12").Value = 1.53
$ws.Range("A13").Value = "ma.L2"
$ws.Range("B13").Value = 0.3765
$ws.Range("C13").Value = 0.115
$ws.Range("D13").Value = 3.283
$ws.Range("E13").Value = 0.001
$ws.Range("F13").Value = 0.152
$ws.Range("G13").Value = 0.601
$ws.Range("A14").Value = "ma.S.L12"
$ws.Range("B14").Value = -0.8548
$ws.Range("C14").Value = 0.083
$ws.Range("D14").Value = -10.282
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = -1.018
$ws.Range("G14").Value = -0.692
$ws.Range("A15").Value = "ma.S.L24"
$ws.Range("B15").Value = 0.2257
$ws.Range("C15").Value = 0.123
$ws.Range("D15").Value = 1.835
$ws.Range("E15").Value = 0.066
$ws.Range("F15").Value = -0.015
$ws.Range("G15").Value = 0.467
$ws.Range("A16").Value = "sigma2"
$ws.Range("B16").Value = 2313000000
$ws.Range("C16").Value = 20.776
$ws.Range("D16").Value = 111000000
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 2310000000
$ws.Range("G16").Value = 2310000000

"Edit complete"
